$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Extension Payments" Tax Type from execution: mark C4 as DONOTRUN
$ws.Range("C4").Value = "DONOTRUN"

# Widen column C to fit the new value
$ws.Range("C1").ColumnWidth = 16

# Update the active selection to C4
$ws.Range("C4").Select()
